# Finalizacion login y registro de usuario
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- New "Rol" table (K59:L63) ---
$ws.Range("K59").Value = "Rol"

$ws.Range("K60").Value = "idRol"
$ws.Range("L60").Value = "nombre"

$ws.Range("K61").Value = 1
$ws.Range("L61").Value = "Lider"

$ws.Range("K62").Value = 2
$ws.Range("L62").Value = "Stakeholder"

$ws.Range("K63").Value = 3
$ws.Range("L63").Value = "Programador"

# Style the new "Rol" table like the other yellow-highlighted tables
# (column K gets the yellow "id" highlight + border; column L only the border)
$ws.Range("K60:L63").Borders.Color = 0
$ws.Range("K60:L63").Borders.LineStyle = 1
$ws.Range("K60:K63").Interior.Color = 65535

# --- New "Rol" column (R) on the "Asignado a" table (rows 43-46) ---
$ws.Range("R43").Value = "Rol"
$ws.Range("R44").Value = 1
$ws.Range("R45").Value = 1
$ws.Range("R46").Value = 3

$ws.Range("R43:R46").Borders.Color = 0
$ws.Range("R43:R46").Borders.LineStyle = 1
$ws.Range("R43:R46").Interior.Color = 65535

# --- sheetView state update ---
$ws.Activate()
$ws.Range("Q51").Select()
